$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old column A (taxon id values 3,5,7,13 with the header style but no header text)
# is removed entirely; the remaining columns B:F shift left to become A:E.
$ws.Range("A:A").EntireColumn.Delete()
